$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# The fork can now hold the mobile goal (mogo) itself, so the second
# pneumatic mogo-lock (row 28 / PORT G) is no longer needed.
# Remove its DEVICE/NAME(code) entries, leaving only the PORT letter.
$ws.Range("B28:C28").ClearContents()

# Update the view so the previously-hidden row 28 area is visible and
# the last touched cell (C28) is selected, matching the saved view state.
$ws.Range("C28").Select()
